# Finalize v1.0 invoice automation release
# Appends new vendor records to the VendorMaster sheet and widens column A
# to fit the newly added (longer) vendor names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vendor rows appended below the existing data (rows 2-23 already used).
# Data entered in the same order it was originally typed in the source
# workbook: names for the first batch of rows, then all of the new vendor
# IDs, then the remaining two vendor names, then the Active flags.
$ws.Cells.Item(24, 1).Value = "Adventure Works Training"
$ws.Cells.Item(25, 1).Value = "Litware Utilities"
$ws.Cells.Item(26, 1).Value = "Northwind Office Supplies"
$ws.Cells.Item(29, 1).Value = "Tailspin Logistics"

$ws.Cells.Item(24, 2).Value = "V023"
$ws.Cells.Item(25, 2).Value = "V024"
$ws.Cells.Item(26, 2).Value = "V025"
$ws.Cells.Item(27, 2).Value = "V026"
$ws.Cells.Item(28, 2).Value = "V027"
$ws.Cells.Item(29, 2).Value = "V028"

$ws.Cells.Item(27, 1).Value = "Contoso IT Services Inc."
$ws.Cells.Item(28, 1).Value = "Fabrikam Facilities Co."

$ws.Cells.Item(24, 3).Value = "Y"
$ws.Cells.Item(25, 3).Value = "Y"
$ws.Cells.Item(26, 3).Value = "N"
$ws.Cells.Item(27, 3).Value = "Y"
$ws.Cells.Item(28, 3).Value = "Y"
$ws.Cells.Item(29, 3).Value = "Y"

# Widen column A to fit the new, longer vendor names (no longer "best fit").
$ws.Columns.Item(1).ColumnWidth = 30.166666666666668

# Clear the stale selection left over on the old last row so the saved
# workbook opens with the default (top-left) selection.
[void]$ws.Range("A1").Select()
